$d = $word.ActiveDocument

$replacements = @(
    @("283×8=2264", "574×7=4018"),
    @("496×9=4464", "459×7=3213"),
    @("568×9=5112", "385×7=2695"),
    @("851×6=5106", "134×8=1072"),
    @("996×5=4980", "445×6=2670"),
    @("843×9=7587", "198×6=1188"),
    @("371×9=3339", "237×7=1659"),
    @("895×8=7160", "822×6=4932"),
    @("693×2=1386", "599×8=4792"),
    @("297×7=2079", "806×6=4836"),
    @("942×3=2826", "740×5=3700"),
    @("796×6=4776", "649×6=3894"),
    @("854×6=5124", "719×6=4314"),
    @("446×3=1338", "303×5=1515"),
    @("736×8=5888", "296×8=2368"),
    @("227×7=1589", "319×7=2233"),
    @("282×5=1410", "696×5=3480"),
    @("708×2=1416", "728×2=1456"),
    @("654×3=1962", "548×8=4384"),
    @("923×8=7384", "343×8=2744"),
    @("922×3=2766", "604×3=1812"),
    @("357×2=714",  "788×7=5516"),
    @("222×6=1332", "894×9=8046"),
    @("481×4=1924", "356×5=1780"),
    @("521×7=3647", "527×7=3689")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
